# Applies the "Resolución" header/value restructuring to vehiculos_prueba.xlsx
#
# Summary of change:
#   C1: "Resolución Padre"      -> "Resolución Primigenia"
#   D1: "Resolución Primigenia" -> "Resolución Hija"
#   D2: "001-2024-DRTC-PUNO"    -> (cleared)
#   D3: "002-2024-DRTC-PUNO"    -> (cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("C1").Value = "Resolución Primigenia"
$ws.Range("D1").Value = "Resolución Hija"

# Clear the now-redundant "Resolución Hija" values on the sample rows;
# column C (Resolución Primigenia) keeps its original value.
$ws.Range("D2").Value = $null
$ws.Range("D3").Value = $null
